$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.658.70'
$ws.Range("E2").Value = '  -7.02%  '

$ws.Range("D3").Value = '1.695.74'

$ws.Range("E4").Value = '  +0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5124'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -12.83%  '

$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("E8").Value = '  -4.60%  '

$ws.Range("E9").Value = '  -4.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06288'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07356'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.31%  '

$ws.Range("D12").Value = '1.697.55'
$ws.Range("E12").Value = '  -5.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.517'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5788'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.62%  '

$ws.Range("D15").Value = '1.926.87'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008437'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.54'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = '26.680.90'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.995'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '186.66'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -11.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.246'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.506'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1156'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.81%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.357'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05649'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.335'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.505'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.491'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.649'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.15%  '

$ws.Range("E35").Value = '  -2.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6006'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.360'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.701'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.49%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01615'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.89%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.102.44'
$ws.Range("E40").Value = '  -3.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8587'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.81%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.841'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.17%  '

$ws.Range("E43").Value = '  -0.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.68'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.57%  '

$ws.Range("D45").Value = '1.853.34'
$ws.Range("E45").Value = '  -5.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000113'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.122'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.45%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05240'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4325'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.50%  '
